# Hana_T1002.xlsx — "Address lookup test data and locators fix done"
#
# The recipient full-address locator string had a typo/format issue:
# "3286 B Hwy 100, Villa Ridge, MO 63089" -> "3286b Hwy 100, Villa Ridge, MO 63089"
# It lives in cell I2 of Sheet1 (header "recipientfulladdress1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = "3286b Hwy 100, Villa Ridge, MO 63089"

# A plain .Value write resets the cell's quote-prefixed-text style; restore
# it by pasting just the number format/style from an unaffected neighbor
# cell that uses the same formatting (border/quotePrefix) as I2 originally did.
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Update the active selection to match the post-edit UI state.
$ws.Range("I5").Select()
